# Add a new comment paragraph after the "GitHub repo..." paragraph, and
# move the _GoBack bookmark (currently sitting mid-sentence in the first
# paragraph) so it trails the newly-added sentence instead.

$d = $word.ActiveDocument

# First paragraph ("Test GitHub repo using a Word doc.") currently holds
# the stray _GoBack bookmark; it will be removed from there once we add
# a fresh _GoBack bookmark elsewhere (Word only keeps one bookmark per
# name, so re-adding it relocates it).
$p1 = $d.Paragraphs(1)

# Insert a brand-new, empty paragraph right after it.
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)

# Fill the new paragraph with the comment text, padding with a throwaway
# trailing character so the bookmark-insertion point below lands strictly
# after the sentence rather than snapping to the paragraph-end position
# (which would otherwise cause it to swallow/merge across the paragraph
# mark).
$p2.Range.Text = "This is a totally new comment.X"

$pEnd = $p2.Range.End

# Bookmark the single placeholder character "X" ...
$bmRange = $d.Range($pEnd - 2, $pEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ... then delete that placeholder character, leaving a zero-length
# bookmark immediately after "This is a totally new comment."
$xRange = $d.Range($pEnd - 2, $pEnd - 1)
$xRange.Text = ""
